# Update the "丽水-漫展信息" workbook:
#   - The 青田 event (row 2) is replaced by the 缙云 event (previously row 3),
#     with an updated "想去人数" (F) value of 54.
#   - The now-duplicate row 3 is removed.
# This applies to both the "展览" sheet and the "全部类型" sheet (sheet1 and
# sheet4), which carry the same two data rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2 (青田 event) is replaced by row 3's (缙云 event) data. Copy/paste
    # (rather than re-typing the date string) keeps the cell a plain text
    # value instead of Excel auto-coercing "2024-10-03" into a date serial.
    $ws.Range("B3:I3").Copy()
    $ws.Range("B2").PasteSpecial()

    # The "想去人数" count was refreshed to 54 (differs from the old row 3's 53).
    $ws.Range("F2").Value = 54

    # Remove the now-redundant row 3, shrinking the used range to A1:I2.
    $ws.Rows.Item(3).Delete()
}
